$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# 1. Merge the split run in the "03.02.2023" row's description cell into a
#    single run containing the full text (the search text is unique in the
#    document, so this only touches that one cell).
$descCell = $table.Rows.Item(4).Cells.Item(3)
$descCell.Range.Find.Execute(
    "Työstön aloituspalaveri ja vaatimusmäärittelyn kirjoitus",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Työstön aloituspalaveri ja vaatimusmäärittelyn kirjoitus", 1
) | Out-Null

# 2. Fill in the next empty row with the new date, hours and description.
$row = $table.Rows.Item(5)

$dateCell = $row.Cells.Item(1)
$dateRng = $dateCell.Range
$dateRng.End = $dateRng.End - 1
$dateRng.Text = "08.02.2023"

$hoursCell = $row.Cells.Item(2)
$hoursRng = $hoursCell.Range
$hoursRng.End = $hoursRng.End - 1
$hoursRng.Text = "1,5h"

# The description cell already holds a single-space placeholder run; replace
# just that occurrence (searching forward from this cell onward so no other
# part of the document is affected) with the real text.
$descCell2 = $row.Cells.Item(3)
$tailRange = $d.Range($descCell2.Range.Start, $d.Content.End)
$tailRange.Find.Execute(
    " ", $true, $false, $false, $false, $false, $true, 0, $false,
    "Luotu alku käyttäjänluonnille", 1
) | Out-Null
